# ontoDog_input.xlsx — "added missing OBI terms and fixed survey data axiom"
#
# Adds 20 new lookup rows (rows 102-121) to Sheet1 with newly-introduced OBI /
# IAO / GAZ / REO term URIs, their labels, and a "Y" flag column. Row 107
# (IAO_0000408 / "length measurement datum") was the fix for the survey data
# axiom and was appended to the sheet after the other 19 rows, so its shared
# strings land at the end of the shared-string table even though the row
# itself sits between rows 106 and 108 in the sheet — we replicate that by
# writing it last while still addressing cell row 107 directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- new data rows -------------------------------------------------------
# (row, colA URL, colB label) -- colC is always the literal "Y" flag.
$ws.Cells.Item(102, 1).Value = 'http://purl.obolibrary.org/obo/GAZ_00000448'
$ws.Cells.Item(102, 2).Value = 'geographic location'
$ws.Cells.Item(102, 3).Value = 'Y'
$ws.Cells.Item(103, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000029'
$ws.Cells.Item(103, 2).Value = 'numeral'
$ws.Cells.Item(103, 3).Value = 'Y'
$ws.Cells.Item(104, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000100'
$ws.Cells.Item(104, 2).Value = 'data set'
$ws.Cells.Item(104, 3).Value = 'Y'
$ws.Cells.Item(105, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000140'
$ws.Cells.Item(105, 2).Value = 'setting datum'
$ws.Cells.Item(105, 3).Value = 'Y'
$ws.Cells.Item(106, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000300'
$ws.Cells.Item(106, 2).Value = 'textual entity'
$ws.Cells.Item(106, 3).Value = 'Y'
$ws.Cells.Item(108, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000422'
$ws.Cells.Item(108, 2).Value = 'postal address'
$ws.Cells.Item(108, 3).Value = 'Y'
$ws.Cells.Item(109, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000429'
$ws.Cells.Item(109, 2).Value = 'email address'
$ws.Cells.Item(109, 3).Value = 'Y'
$ws.Cells.Item(110, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000577'
$ws.Cells.Item(110, 2).Value = 'CRID symbol'
$ws.Cells.Item(110, 3).Value = 'Y'
$ws.Cells.Item(111, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000578'
$ws.Cells.Item(111, 2).Value = 'CRID'
$ws.Cells.Item(111, 3).Value = 'Y'
$ws.Cells.Item(112, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000582'
$ws.Cells.Item(112, 2).Value = 'time stamped measurement datum'
$ws.Cells.Item(112, 3).Value = 'Y'
$ws.Cells.Item(113, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0000102'
$ws.Cells.Item(113, 2).Value = 'responsible party role'
$ws.Cells.Item(113, 3).Value = 'Y'
$ws.Cells.Item(114, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0000103'
$ws.Cells.Item(114, 2).Value = 'principal investigator role'
$ws.Cells.Item(114, 3).Value = 'Y'
$ws.Cells.Item(115, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0000274'
$ws.Cells.Item(115, 2).Value = 'adding a material entity into a target'
$ws.Cells.Item(115, 3).Value = 'Y'
$ws.Cells.Item(116, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0000444'
$ws.Cells.Item(116, 2).Value = 'target of material addition role'
$ws.Cells.Item(116, 3).Value = 'Y'
$ws.Cells.Item(117, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0001936'
$ws.Cells.Item(117, 2).Value = 'molecular-labeled material'
$ws.Cells.Item(117, 3).Value = 'Y'
$ws.Cells.Item(118, 1).Value = 'http://purl.obolibrary.org/obo/OBI_0302914'
$ws.Cells.Item(118, 2).Value = 'digital curation'
$ws.Cells.Item(118, 3).Value = 'Y'
$ws.Cells.Item(119, 1).Value = 'http://purl.obolibrary.org/obo/OBI_1110087'
$ws.Cells.Item(119, 2).Value = 'donor role'
$ws.Cells.Item(119, 3).Value = 'Y'
$ws.Cells.Item(120, 1).Value = 'http://purl.obolibrary.org/obo/OBI_1110109'
$ws.Cells.Item(120, 2).Value = 'target of material addition'
$ws.Cells.Item(120, 3).Value = 'Y'
$ws.Cells.Item(121, 1).Value = 'http://purl.obolibrary.org/obo/REO_0000280'
$ws.Cells.Item(121, 2).Value = 'molecular label'
$ws.Cells.Item(121, 3).Value = 'Y'

# Row 107 (the survey-data-axiom fix) was appended last, after the batch of
# 19 other rows above, even though it is sandwiched between rows 106/108 in
# the sheet -- so its own shared strings trail the rest in sharedStrings.xml.
$ws.Cells.Item(107, 1).Value = 'http://purl.obolibrary.org/obo/IAO_0000408'
$ws.Cells.Item(107, 2).Value = 'length measurement datum'
$ws.Cells.Item(107, 3).Value = 'Y'

# --- window / view bookkeeping (best effort) ------------------------------
$win = $excel.ActiveWindow
$win.Left = 3460
$win.Top = 3780
$win.ScrollRow = 80
$win.ScrollColumn = 1
$ws.Range("B107").Select()
